# Weekly update: insert a new price record as row 89, shifting existing
# rows 89:173 down to 90:174 (matches commit "Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("89:89").Insert()

$ws.Range("A89").Value = 7
$ws.Range("B89").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C89").Value = "Ñuble"
$ws.Range("D89").Value = 44762
$ws.Range("E89").Value = 16
$ws.Range("F89").Value = "Fruta"
$ws.Range("G89").Value = 100102
$ws.Range("H89").Value = "Cítricos"
$ws.Range("I89").Value = 100102004
$ws.Range("J89").Value = "Mandarina"
$ws.Range("K89").Value = "Clemenuless"
$ws.Range("L89").Value = "Primera"
$ws.Range("M89").Value = 120
$ws.Range("N89").Value = 7500
$ws.Range("O89").Value = 8000
$ws.Range("P89").Value = 7750
$ws.Range("Q89").Value = "`$/caja 18 kilos"
$ws.Range("R89").Value = "Región de O'Higgins"
$ws.Range("S89").Value = 431
$ws.Range("T89").Value = 18
